# SPIDAL-Tutorial-Feb2017-midas.pptx edit
# The MDAnalysis SPIDAL tutorial moved to a new GitHub org / doc URL.
# Update the two affected hyperlink texts on slide 1 and refresh the
# handout-master date placeholder.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item(1)
$tr = $shape.TextFrame.TextRange
$full = $tr.Text

# --- 1) "https://github.com/radical-cybertools/radical.pilot/tree" + "/master"
#        -> single run "https://github.com/radical-cybertools/radical.pilot/tree/master"
$oldUrl1 = "https://github.com/radical-cybertools/radical.pilot/tree/master"
$idx1 = $full.IndexOf($oldUrl1)
if ($idx1 -ge 0) {
    $run1 = $tr.Characters($idx1 + 1, $oldUrl1.Length)
    $run1.Text = $oldUrl1
}

# --- 2) MDAnalysis tutorial link moved from
#        https://becksteinlab.github.io/SPIDAL-MDAnalysis-Midas-tutorial/index.html
#        to
#        http://www.mdanalysis.org/SPIDAL-MDAnalysis-Midas-tutorial/index.html
$full = $tr.Text
$oldUrl2 = "https://becksteinlab.github.io/SPIDAL-MDAnalysis-Midas-tutorial/index.html"
$idx2 = $full.IndexOf($oldUrl2)
if ($idx2 -ge 0) {
    $start = $idx2 + 1

    # Replace the whole old URL with the new one first (collapses it to a
    # single run using the first run's formatting/hyperlink).
    $newUrl2 = "http://www.mdanalysis.org/SPIDAL-MDAnalysis-Midas-tutorial/index.html"
    $whole = $tr.Characters($start, $oldUrl2.Length)
    $whole.Text = $newUrl2

    # Re-split the new URL into the four runs seen in the authored deck:
    #   "http:"  "/"  "/www.mdanalysis.org/"  "SPIDAL-MDAnalysis-Midas-tutorial/index.html"
    $part1 = "http:"
    $part2 = "/"
    $part3 = "/www.mdanalysis.org/"
    $part4 = "SPIDAL-MDAnalysis-Midas-tutorial/index.html"

    $r1 = $tr.Characters($start, $part1.Length)
    $r1.Text = $part1

    $r2 = $tr.Characters($start + $part1.Length, $part2.Length)
    $r2.Text = $part2

    $r3 = $tr.Characters($start + $part1.Length + $part2.Length, $part3.Length)
    $r3.Text = $part3

    $r4 = $tr.Characters($start + $part1.Length + $part2.Length + $part3.Length, $part4.Length)
    $r4.Text = $part4
}

# --- 3) Handout master "date" placeholder text 2/8/2017 -> 2017-02-21
$hm = $p.HandoutMaster
$hf = $hm.HeadersFooters
$hf.DateAndTime.Text = "2017-02-21"
